# Applies the "material-1 data & procedure of material-2" edit:
#   1. Title line "Exp1 ..." - normalize run split (E|xp|1|' '|title -> Exp|1|' '|title)
#      Net visible text is unchanged; this just re-types the existing text so the
#      redundant run boundary collapses the way it does in the target document.
#   2. Author list "主试编号：E1-伍嘉琪 " gains two more experimenters:
#      "主试编号：E1-伍嘉琪、E2-孙心茹、E3-张绍明"
#   3. The two pilot-study hyperlink URLs ("pilot_friend/" and "pilot_stranger/")
#      were previously split mid-word across three runs each; re-type them as a
#      single contiguous run (same visible URL/text, same hyperlink target).

$d = $word.ActiveDocument

# --- 1. Title run normalization -------------------------------------------------
$d.Content.Find.Execute("Exp1 ", $true, $false, $false, $false, $false, $true, `
    [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue, $false, "Exp1 ", `
    [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceAll) | Out-Null

# --- 2. Add E2 / E3 experimenters after E1-伍嘉琪 --------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("主试编号：E1-伍嘉琪", $true, $false, $false, $false, $false, $true, `
    [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue, $false, "", `
    [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceNone)
if ($found) {
    $rng.Collapse([Microsoft.Office.Interop.Word.WdCollapseDirection]::wdCollapseEnd) | Out-Null
    $rng.MoveEnd([Microsoft.Office.Interop.Word.WdUnits]::wdCharacter, 1) | Out-Null
    $rng.Text = "、E2-孙心茹、E3-张绍明"
}

# --- 3. Merge the split hyperlink URL runs --------------------------------------
$d.Content.Find.Execute("pilot_friend/", $true, $false, $false, $false, $false, $true, `
    [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue, $false, "pilot_friend/", `
    [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceAll) | Out-Null

$d.Content.Find.Execute("pilot_stranger/", $true, $false, $false, $false, $false, $true, `
    [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue, $false, "pilot_stranger/", `
    [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceAll) | Out-Null
